# edit.ps1
# -----------------------------------------------------------------------
# Reproduces the commit:
#   1. The table on slide 5 (the 3-column "types of financial documents"
#      table) gets re-styled: its tableStyleId changes from the pink
#      "Table_0" custom style to the built-in style
#      {9401E8EC-4985-403E-BBCE-7CD5F4CD1D3D}.
#   2. The deck's theme colours (ppt/theme/theme1.xml, which drives the
#      slide master/every slide) are swapped from the "Integral" /
#      "Red Violet" palette to the stock "Office Theme" palette.
#      (dk1/lt1/fonts/format-scheme are already identical between the
#      two palettes, so only the 10 colours that actually differ need
#      to be touched.)
# -----------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 ------------------------------------------------
$tableSlide  = $p.Slides.Item(5)
$tableShape  = $tableSlide.Shapes.Item(2)          # the graphicFrame holding the table
$table       = $tableShape.Table
$table.ApplyStyle("{9401E8EC-4985-403E-BBCE-7CD5F4CD1D3D}", $true)

# --- 2. Theme colours (Integral/Red Violet -> Office Theme) ------------------
# Use the first slide purely as a handle onto the shared ThemeColorScheme -
# every slide shares the same slide master / theme part, so this updates
# ppt/theme/theme1.xml for the whole deck.
$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme

# Index -> theme slot: 1=dk1 2=lt1 3=dk2 4=lt2 5-10=accent1-6 11=hlink 12=folHlink
# dk1 (000000) and lt1 (FFFFFF) are unchanged between the two palettes.
$themeColors.Colors(3).RGB  = 6968388    # dk2      44546A
$themeColors.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$themeColors.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$themeColors.Colors(6).RGB  = 3243501    # accent2  ED7D31
$themeColors.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$themeColors.Colors(8).RGB  = 49407      # accent4  FFC000
$themeColors.Colors(9).RGB  = 12874308   # accent5  4472C4
$themeColors.Colors(10).RGB = 4697456    # accent6  70AD47
$themeColors.Colors(11).RGB = 12673797   # hlink    0563C1
$themeColors.Colors(12).RGB = 7491477    # folHlink 954F72

# Best-effort: also try to rename the design/theme to match ("Office
# Theme"); harmless if the host doesn't support renaming a Design.
try {
    $design = $p.Designs.Item(1)
    $design.Name = "Office Theme"
} catch {
}
